# Automatische test-sync: 2025-06-24 21:36:50
# Append the new mail-log entry to the "Logs" sheet, extend the
# conditional-formatting ranges to cover it, and refresh the
# "Dashboard" summary sheet to reflect the updated category counts.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append new row to the Logs sheet (row 34) ---
$newRow = 34

$logs.Cells.Item($newRow, 1).Value = "Order wijzigen"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Kan ik mijn bestelling nog aanpassen?"
$logs.Cells.Item($newRow, 4).Value = "Bestelling / Levering"
$logs.Cells.Item($newRow, 5).Value = "nan"
$logs.Cells.Item($newRow, 6).Value = "2025-06-24 21:36:43"
$logs.Cells.Item($newRow, 7).Value = "Ja"

# --- Extend the conditional formatting ranges to include the new row ---
$dFormatConditions = $logs.Range("D2:D33").FormatConditions
for ($i = 1; $i -le $dFormatConditions.Count; $i++) {
    $dFormatConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D34"))
}

$gFormatConditions = $logs.Range("G2:G33").FormatConditions
for ($i = 1; $i -le $gFormatConditions.Count; $i++) {
    $gFormatConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G34"))
}

# --- Update the Dashboard summary table ---
# "Bestelling / Levering" count increases to 4 and now ranks above
# "Offerte / Prijsaanvraag" (count stays 3) in the sorted summary.
$dashboard.Cells.Item(4, 1).Value = "Bestelling / Levering"
$dashboard.Cells.Item(4, 2).Value = 4

$dashboard.Cells.Item(6, 1).Value = "Offerte / Prijsaanvraag"
$dashboard.Cells.Item(6, 2).Value = 3
